$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 131916.49741240995
$ws.Range("C2").Value = 234.61669921875
$ws.Range("D2").Value = 235.08706665039062
$ws.Range("E2").Value = 236.3173370361328
$ws.Range("F2").Value = 9.397096633911133
$ws.Range("G2").Value = 8.64663314819336
$ws.Range("H2").Value = 7.119923114776611
$ws.Range("I2").Value = 1984.673095703125
$ws.Range("J2").Value = 1562.0362548828125
$ws.Range("K2").Value = 1612.2728271484375
$ws.Range("L2").Value = 0.8556978106498718
$ws.Range("M2").Value = 0.7689777612686157
$ws.Range("N2").Value = 0.9590601921081543
$ws.Range("O2").Value = 5158.982421875

# Row 3
$ws.Range("B3").Value = 28434.061172639904
$ws.Range("C3").Value = 234.65138244628906
$ws.Range("D3").Value = 234.8115234375
$ws.Range("E3").Value = 236.37496948242188
$ws.Range("F3").Value = 1.0488389730453491
$ws.Range("G3").Value = 1.7748174667358398
$ws.Range("H3").Value = 1.5231413841247559
$ws.Range("I3").Value = 226.21859741210938
$ws.Range("J3").Value = 259.76190185546875
$ws.Range("K3").Value = 315.5941467285156
$ws.Range("L3").Value = 0.919171154499054
$ws.Range("M3").Value = 0.6233075261116028
$ws.Range("N3").Value = 0.8765712976455688
$ws.Range("O3").Value = 801.5746459960938

# Row 4
$ws.Range("B4").Value = 168542.32664242078
$ws.Range("C4").Value = 234.83279418945312
$ws.Range("D4").Value = 235.4658660888672
$ws.Range("E4").Value = 236.6804962158203
$ws.Range("F4").Value = 0.2010851800441742
$ws.Range("G4").Value = 0.22693611681461334
$ws.Range("H4").Value = 0.45828041434288025
$ws.Range("I4").Value = 31.98046112060547
$ws.Range("J4").Value = -17.79543113708496
$ws.Range("K4").Value = 51.478572845458984
$ws.Range("L4").Value = 0.6772451996803284
$ws.Range("M4").Value = 0.33302509784698486
$ws.Range("N4").Value = 0.4746054708957672
$ws.Range("O4").Value = 65.80813598632812

# Row 5
$ws.Range("B5").Value = 126109.6145374211
$ws.Range("C5").Value = 235.8601531982422
$ws.Range("D5").Value = 234.41244506835938
$ws.Range("E5").Value = 234.38795471191406
$ws.Range("F5").Value = 56.75860595703125
$ws.Range("G5").Value = 53.0229606628418
$ws.Range("H5").Value = 55.24835968017578
$ws.Range("I5").Value = 10410.4560546875
$ws.Range("J5").Value = 9372.2060546875
$ws.Range("K5").Value = 10191.412109375
$ws.Range("L5").Value = 0.7776487469673157
$ws.Range("M5").Value = 0.7540448904037476
$ws.Range("N5").Value = 0.7870090007781982
$ws.Range("O5").Value = 29974.07421875

# Row 6
$ws.Range("B6").Value = 26861.310665731053
$ws.Range("C6").Value = 236.27127075195312
$ws.Range("D6").Value = 234.56773376464844
$ws.Range("E6").Value = 234.64788818359375
$ws.Range("F6").Value = 0.2897227108478546
$ws.Range("G6").Value = 0.40349316596984863
$ws.Range("H6").Value = 0.2573781907558441
$ws.Range("I6").Value = -19.418872833251953
$ws.Range("J6").Value = 83.25874328613281
$ws.Range("K6").Value = 10.165909767150879
$ws.Range("L6").Value = 0.28375619649887085
$ws.Range("M6").Value = 0.8797476887702942
$ws.Range("N6").Value = 0.16836977005004883
$ws.Range("O6").Value = 74.00578308105469

# Row 7
$ws.Range("B7").Value = 24350.219579162927
$ws.Range("C7").Value = 234.1078643798828
$ws.Range("D7").Value = 235.02101135253906
$ws.Range("E7").Value = 236.3190155029297
$ws.Range("F7").Value = 13.248489379882812
$ws.Range("G7").Value = 4.182949066162109
$ws.Range("H7").Value = 3.9002068042755127
$ws.Range("I7").Value = 2986.270751953125
$ws.Range("J7").Value = 571.7822875976562
$ws.Range("K7").Value = 479.17529296875
$ws.Range("L7").Value = 0.962823748588562
$ws.Range("M7").Value = 0.5810154676437378
$ws.Range("N7").Value = 0.5213376879692078
$ws.Range("O7").Value = 4049.338134765625

# Row 8
$ws.Range("B8").Value = 116594.76235985631
$ws.Range("C8").Value = 234.5697021484375
$ws.Range("D8").Value = 234.69024658203125
$ws.Range("E8").Value = 236.71922302246094
$ws.Range("F8").Value = 26.598766326904297
$ws.Range("G8").Value = 35.961387634277344
$ws.Range("H8").Value = 32.13772201538086
$ws.Range("I8").Value = 5790.56689453125
$ws.Range("J8").Value = 7175.96533203125
$ws.Range("K8").Value = 7237.93505859375
$ws.Range("L8").Value = 0.9280848503112793
$ws.Range("M8").Value = 0.8502544164657593
$ws.Range("N8").Value = 0.9514064192771912
$ws.Range("O8").Value = 20204.466796875
